$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 4672
$ws1.Range("F3").Value = 1853
$ws1.Range("F9").Value = 275
$ws1.Range("F11").Value = 543
$ws1.Range("F12").Value = 537
$ws1.Range("F13").Value = 390
$ws1.Range("F17").Value = 127
$ws1.Range("F22").Value = 12
$ws1.Range("F26").Value = 53
$ws1.Range("F32").Value = 3892
$ws1.Range("F34").Value = 769
$ws1.Range("F36").Value = 1073
$ws1.Range("F38").Value = 1859

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 4672
$ws4.Range("F3").Value = 1853
$ws4.Range("F9").Value = 275
$ws4.Range("F11").Value = 543
$ws4.Range("F12").Value = 537
$ws4.Range("F14").Value = 390
$ws4.Range("F18").Value = 127
$ws4.Range("F23").Value = 12
$ws4.Range("F27").Value = 53
$ws4.Range("F33").Value = 3892
$ws4.Range("F36").Value = 769
$ws4.Range("F38").Value = 1073
$ws4.Range("F40").Value = 1859
